$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notifications")

# --- Update "Send To" text for the On Submit Challenge / On First Vote rows ---
$ws.Range("B2").Value = "Donator - Creator - Admin"
$ws.Range("B3").Value = "Challenger - Admin"

# --- Add new "Web" status column (G) for every row that already has Done/Api ---
$ws.Range("G2").Value = "Web"
$ws.Range("G3").Value = "Web"
$ws.Range("G5").Value = "Web"
$ws.Range("G6").Value = "Web"
$ws.Range("G7").Value = "Web"
$ws.Range("G8").Value = "Web"

# Match alignment used by the rest of the data cells for the new column
$ws.Range("G2:G8").HorizontalAlignment = -4108
$ws.Range("G2:G8").VerticalAlignment = -4108

# --- Swap row 9 ("On Create Challage") and row 10 ("On Win") and add Web to the win row ---
$ws.Range("A9").Value = "On Win "
$ws.Range("B9").Value = "To Donator, Creator, Submitor ,Winner & Admin"
$ws.Range("C9").Value = "SUBMITED_CHALLENGE_LIST_SCREEN"
$ws.Range("D9").Value = "Challenge ID"
$ws.Range("E9").Value = "Done"
$ws.Range("F9").Value = "Api"
$ws.Range("G9").Value = "Web"
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G9").VerticalAlignment = -4108
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108

$ws.Range("A10").Value = "On Create Challage"
$ws.Range("B10").Value = "to Admin"
$ws.Range("C10").Value = "ADMIN_NOTIFICATION"
$ws.Range("D10").Value = "Challenge ID"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""

# --- Drop the now-unused trailing blank row (was row 23) ---
$ws.Range("A23:G23").Value = ""

$ws.Range("G10").Select()
